# "Generate Report for handoff"
#
# The handoff run failed to produce a transform, so the status report is
# regenerated to reflect that: the per-language "Latest Handoff File"
# hyperlink/cell is cleared (no file was produced), the handoff/handback
# datetimes reset to the zero value, the handoff reason becomes "Ignored",
# and the overall status changes from "Ready for handoff" to
# "Handoff transform failed" everywhere it is shown.

$wb = $excel.ActiveWorkbook

$zeroDate = "0001-01-01 00:00:00"
$newStatus = "Handoff transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Update the "Status" text everywhere it is referenced -------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsZhCn.Range("B2").Value = $newStatus
$wsDeDe.Range("B2").Value = $newStatus

# --- Per-language detail sheets ---------------------------------------
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    # Drop the "Latest Handoff File" hyperlink + cell content for row 2 -
    # no handoff file was produced by the failed transform.
    $toRemove = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$C$2') {
            $toRemove += $h
        }
    }
    foreach ($h in $toRemove) {
        $h.Delete()
    }
    $ws.Range("C2").Clear()

    # Reset the handoff/handback datetimes to the zero value.
    $ws.Range("D2").Value = $zeroDate
    $ws.Range("G2").Value = $zeroDate
    $ws.Range("D3").Value = $zeroDate
    $ws.Range("G3").Value = $zeroDate

    # The dependency is now ignored rather than included.
    $ws.Range("H2").Value = "Ignored"
    $ws.Range("H3").Value = "Ignored"
}
